$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ast_config")

$ws.Range("M2").Value = "Queued"
$ws.Range("M3").Value = "Queued"
